# Reorder the "BrandItemPageData" sheet so it comes after "BrandPageData"
# (this is the "item crawling action" sheet being re-positioned/recreated).
$wb = $excel.ActiveWorkbook

$itemSheet = $wb.Worksheets.Item("BrandItemPageData")
$brandSheet = $wb.Worksheets.Item("BrandPageData")
$itemSheet.Move($null, $brandSheet)

# Update the breadcrumb text cells to reflect the new crawling action:
# "# Famous" -> "# Famous Cigars for Sale"
$ws = $wb.Worksheets.Item("BrandItemPageData")
$ws.Range("G2").Value = "» Famous Smoke Shop Cigars » Discount Cigars Online » Cigar Brand List » # Famous Cigars for Sale » # Famous Petite Corona Cigars - Natural"
$ws.Range("G3").Value = "» Famous Smoke Shop Cigars » Discount Cigars Online » Cigar Brand List » # Famous Cigars for Sale » # Famous Toro Cigars - Natural"
